$wb = $excel.ActiveWorkbook

# Update both the "展览" and "全部类型" sheets, which contain the same
# table of data. Cell F3 (想去人数 / "want to go" count) increases from
# 75 to 76, F4 from 1459 to 1461, and F9 from 244 to 245.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 76
    $ws.Range("F4").Value = 1461
    $ws.Range("F9").Value = 245
}
